$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ES")
$ws.Range("A31:D31").Copy()
$ws.Range("A32:D32").PasteSpecial(-4122)
$ws.Rows.Item(32).RowHeight = 13.5
$ws.Rows.Item(32).UseStandardHeight = 1
